# Updates cryptos list with new price/volume figures (and a row26/row27 swap)
# matching the commit 'Updated cryptos list ... with GitHub Actions'.
# Numeric-looking Price values are written with a leading apostrophe so Excel
# stores them as text (preserving trailing zeros / exact formatting) instead
# of silently coercing them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.017.60"
$ws.Range("E2").Value = "  +2.39%  "

$ws.Range("D3").Value = "3.188.23"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'536.56"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").Value = "'145.20"
$ws.Range("E6").Value = "  +4.17%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("D9").Value = "'7.35"
$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").Value = "'0.113"
$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("D11").Value = "'0.432"
$ws.Range("E11").Value = "  -1.54%  "

$ws.Range("D12").Value = "3.739.03"
$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("D13").Value = "'0.137"
$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("D14").Value = "'25.81"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").Value = "'0.0000173"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").Value = "60.013.53"
$ws.Range("E16").Value = "  +2.30%  "

$ws.Range("D17").Value = "3.212.79"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("D18").Value = "'6.23"
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").Value = "'13.24"
$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").Value = "'8.19"
$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").Value = "'369.37"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "'0.523"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").Value = "'69.43"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  +1.11%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'8.56"
$ws.Range("E26").Value = "  +3.23%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'0.982"
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("D28").Value = "0.0₃0875"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("D29").Value = "'22.49"
$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "'6.11"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").Value = "'5.27"
$ws.Range("E32").Value = "  +2.98%  "

$ws.Range("E33").Value = "  +5.06%  "

$ws.Range("E34").Value = "  +2.99%  "

$ws.Range("D35").Value = "'157.77"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  +1.71%  "

$ws.Range("D37").Value = "'26.44"
$ws.Range("E37").Value = "  +5.97%  "

$ws.Range("D38").Value = "2.786.16"
$ws.Range("E38").Value = "  +5.93%  "

$ws.Range("D39").Value = "'0.0707"
$ws.Range("E39").Value = "  +3.41%  "

$ws.Range("E40").Value = "  +8.26%  "

$ws.Range("D41").Value = "'1.69"
$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("D42").Value = "'4.22"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").Value = "'39.79"
$ws.Range("E43").Value = "  +2.08%  "

$ws.Range("D44").Value = "'0.719"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "'0.105"
$ws.Range("E45").Value = "  +1.89%  "

$ws.Range("D46").Value = "3.229.39"
$ws.Range("E46").Value = "  +1.21%  "

$ws.Range("D47").Value = "'0.984"
$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("D48").Value = "'6.15"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("D49").Value = "'20.61"
$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("D50").Value = "'0.797"
$ws.Range("E50").Value = "  +6.04%  "

$ws.Range("E51").Value = "  +0.03%  "
